{"js": "// Remove the \"_GoBack\" bookmark (bookmarkStart/bookmarkEnd pair) from the document.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Load all paragraphs so we can find the \">>>  your stuff after this line >>>\" paragraph.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\nlet targetPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"your stuff after this line\") !== -1) {\n    targetPara = paragraphs.items[i];\n    break;\n  }\n}\n\nif (targetPara) {\n  // Re-write the paragraph's whole text as a single run, collapsing the\n  // multiple runs / proofErr markers that used to split it up.\n  const wholeRange = targetPara.getRange(\"Whole\");\n  wholeRange.insertText(\">>>  your stuff after this line >>>\", \"Replace\");\n  await context.sync();\n\n  // Insert a new paragraph right after it containing the red \"This File is\n  // changed.\" text.\n  const newPara = targetPara.insertParagraph(\"This File is changed.\", \"After\");\n  newPara.font.color = \"#FF0000\";\n  await context.sync();\n\n  // Also color the paragraph mark itself (end of paragraph) so the pPr/rPr\n  // carries the same red color as the run.\n  const endRange = newPara.getRange(\"End\");\n  endRange.font.color = \"#FF0000\";\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove the \"_GoBack\" bookmark (bookmarkStart/bookmarkEnd pair).\n$bm = $d.Bookmarks(\"_GoBack\")\n$bm.Delete()\n\n# Collapse the \">>>  your stuff after this line >>>\" paragraph (currently\n# split across several runs with proofErr markers) into a single run by\n# finding its text and replacing it with itself.\n$find = $d.Content.Find\n$find.Text = \">>>  your stuff after this line >>>\"\n$find.Execute(\">>>  your stuff after this line >>>\", $false, $false, $false, $false, $false, $true, 1, $false, \">>>  your stuff after this line >>>\", 2)\n\n# Locate that paragraph again and insert a new paragraph right after it with\n# the red \"This File is changed.\" text.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*your stuff after this line*\") {\n        $p.Range.InsertParagraphAfter()\n        $newPara = $d.Paragraphs.Item($i + 1)\n        $newRange = $newPara.Range\n        $newRange.Text = \"This File is changed.\"\n        $newRange.Font.Color = 255\n        break\n    }\n}\n"}
